# Apply updated "dSF" (column F) values for the marte_yunior 2023 save data.
# This reflects a repull/recalculation of data for several rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -6
    5  = -4
    6  = -7
    7  = 4
    8  = -1
    10 = 0
    18 = -2
    23 = -2
    24 = -5
    27 = 0
    36 = -9
    40 = -7
    41 = -4
    51 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
